$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.876.14'
$ws.Range('E2').Value = '  +2.77%  '
$ws.Range('D3').Value = '1.668.63'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '214.72'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  -0.16%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '23.56'
$c.ClearFormats()
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.904.43'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '1.670.59'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.15'
$c.ClearFormats()
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('E15').Value = '  -1.11%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '66.04'
$c.ClearFormats()
$ws.Range('E16').Value = '  -0.65%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '251.19'
$c.ClearFormats()
$ws.Range('E17').Value = '  +7.11%  '
$ws.Range('D18').Value = '27.834.46'
$ws.Range('E18').Value = '  +2.68%  '
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('E21').Value = '  -0.16%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.47'
$c.ClearFormats()
$ws.Range('E22').Value = '  -1.29%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.34'
$c.ClearFormats()
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('E24').Value = '  -1.46%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '146.91'
$c.ClearFormats()
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('E26').Value = '  -2.91%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '16.32'
$c.ClearFormats()
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +5.85%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0501'
$c.ClearFormats()
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('E33').Value = '  -2.48%  '
$ws.Range('D34').Value = '1.425.15'
$ws.Range('E34').Value = '  -7.36%  '
$ws.Range('E35').Value = '  -5.66%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  -1.33%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.583'
$c.ClearFormats()
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.23'
$c.ClearFormats()
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '5.39'
$c.ClearFormats()
$ws.Range('E44').Value = '  -6.56%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.812.32'
$ws.Range('E45').Value = '  -0.88%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.788'
$c.ClearFormats()
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('E47').Value = '  +4.95%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '89.07'
$c.ClearFormats()
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('E50').Value = '  -1.86%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.84'
$c.ClearFormats()
$ws.Range('E51').Value = '  -4.12%  '
